# Update "want to go" counts (column F) on three sheets, per the commit's
# regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1047
$ws1.Range("F10").Value = 15
$ws1.Range("F13").Value = 90
$ws1.Range("F16").Value = 194
$ws1.Range("F17").Value = 111
$ws1.Range("F21").Value = 178
$ws1.Range("F27").Value = 888

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 27
$ws2.Range("F12").Value = 25

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1047
$ws4.Range("F12").Value = 15
$ws4.Range("F15").Value = 90
$ws4.Range("F18").Value = 194
$ws4.Range("F19").Value = 111
$ws4.Range("F26").Value = 27
$ws4.Range("F29").Value = 178
$ws4.Range("F35").Value = 888
$ws4.Range("F47").Value = 25
